# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD, AE, AF with the same style as the existing
# header cells (e.g. AC1) so they match the bold/centered/bordered look.
$headerStyle = $ws.Range("AC1")
$headerStyle.Copy($ws.Range("AD1"))
$headerStyle.Copy($ws.Range("AE1"))
$headerStyle.Copy($ws.Range("AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-48: every row gets the same season record values.
$lastRow = 48
$wins = 97
$losses = 65
$ties = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
